# Adds a new test-case row (row 29) to the "TestedeSistema" worksheet.
#
# The new row documents one more reservation test case (medico na
# posição 1, sala na posição 7, 2020/11/14 11:00-21:00) that currently
# FAILS ("Fracasso") together with the note describing the fix that is
# still needed, mirroring the layout used by every other row in the
# table (merged C:D / E:F / G:H / I:J / K:L blocks).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 29

# --- Cell values -------------------------------------------------------
$ws.Range("C$row").Value = "O usuario deve fazer uma reserva de sala, inserindo o numero 6, depois escolher o medico na posição 1 na lista, depois a sala de posição 7 na lista, e por fim, escolher a data inicial como 2020/11/14/11/00 e final como 2020/11/14/21/00"
$ws.Range("E$row").Value = "O programa deve impedir que o usuario faça a reserva"
$ws.Range("G$row").Value = "O programa permite que o usuario faça a reserva"
$ws.Range("I$row").Value = "Fracasso"
$ws.Range("K$row").Value = "Ápos concertar o bug, o programa impede que um medico faça uma reserva numa data que colida com a data de uma reserva que esse mesmo medico fez numa sala diferente "

# --- Merge the same column pairs used by every other row ---------------
$ws.Range("C$row`:D$row").MergeCells = $true
$ws.Range("E$row`:F$row").MergeCells = $true
$ws.Range("G$row`:H$row").MergeCells = $true
$ws.Range("I$row`:J$row").MergeCells = $true
$ws.Range("K$row`:L$row").MergeCells = $true

# --- Styling: mirror the rest of the table ------------------------------
# Descrição / Resultado esperado / Resultado obtido / Após revisão
# columns: centered, wrapped (same look as every other data row).
$ws.Range("C$row`:H$row").HorizontalAlignment = -4108  # xlCenter
$ws.Range("C$row`:H$row").VerticalAlignment = -4108    # xlCenter
$ws.Range("C$row`:H$row").WrapText = $true

$ws.Range("K$row`:L$row").HorizontalAlignment = -4108
$ws.Range("K$row`:L$row").VerticalAlignment = -4108
$ws.Range("K$row`:L$row").WrapText = $true

# Conclusão column: centered, but not wrapped (same as I5:J5 "Sucesso").
$ws.Range("I$row`:J$row").HorizontalAlignment = -4108
$ws.Range("I$row`:J$row").VerticalAlignment = -4108

# Every other row uses Excel's maximum row height for its long text.
$ws.Rows.Item($row).RowHeight = 409.5

# --- Update the stored selection / active cell --------------------------
$ws.Range("S$row").Select()
